$d = $word.ActiveDocument

# The document contains a short, standalone paragraph holding only the
# italicized book title "Exodus" directly under the "EXO" Heading2
# paragraph. That whole paragraph (its runs and its paragraph mark) is
# being removed; the trailing empty run that followed "Exodus" becomes
# the new trailing run of the "EXO" heading paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13)
    if ($text -eq "Exodus" -and $p.Range.Italic) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
